$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates as Excel serial numbers, matching existing data)
$data = @(
    @{ Row = 252; A = 44326; B = 9; C = 56; D = 139.0337156760515 },
    @{ Row = 253; A = 44327; B = 4; C = 55; D = 136.5509707532648 },
    @{ Row = 254; A = 44328; B = 0; C = 49; D = 121.654501216545 },
    @{ Row = 255; A = 44329; B = 6; C = 41; D = 101.7925418342519 }
)

foreach ($item in $data) {
    $r = $item.Row

    # Column A: date serial value, formatted/styled like the existing date column
    $ws.Range("A250").Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r, 1).Value2 = $item.A

    # Columns B, C, D: plain numeric values
    $ws.Cells.Item($r, 2).Value2 = $item.B
    $ws.Cells.Item($r, 3).Value2 = $item.C
    $ws.Cells.Item($r, 4).Value2 = $item.D
}
